# Add a new "StatQuery" column to the breed-filter workbook (commit: "updated
# canine breed xl files with stat bar query").
#
# Before:  A=query        B=dbExcel                         C=WebExcel
#          A2=<big query> B2=<Neo4jData.xlsx file name>      C2=<WebData.xlsx file name>
#
# After:   A=query        B=StatQuery                       C=dbExcel                         D=WebExcel
#          A2=<big query> B2=<new stat/count query>          C2=<Neo4jData.xlsx file name>     D2=<WebData.xlsx file name>
#
# i.e. a new column is inserted at B, shifting the old B/C to C/D, and the
# new column is populated with a header label and a Cypher query that
# returns counts (files/samples/cases/studies) instead of row-level data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing column B ("dbExcel"/Neo4j file),
# pushing it (and the WebExcel column after it) one column to the right.
$ws.Columns("B:B").Insert()

# New header label for the inserted column.
$ws.Range("B1").Value = "StatQuery"

# New stat/count Cypher query (counts of files / samples / cases / studies)
# for the West Highland White Terrier breed filter.
$ws.Range("B2").Value = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['West Highland White Terrier']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

# Match the wide "query text" look of column A for the new column (wrap
# text is already inherited from the Insert, this just makes the intent
# explicit and sets the same width as column A).
$ws.Range("B2").WrapText = $true
$ws.Columns("B:B").ColumnWidth = $ws.Columns("A:A").ColumnWidth
